$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 42 (shifts former rows 42-69 down to 44-71)
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(43).Insert()

# New row 42: Jengibre, Primera, fecha 2021-11-29 (serial 44529)
$ws.Cells.Item(42, 1).Value = 9
$ws.Cells.Item(42, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(42, 3).Value = "Metropolitana"
$ws.Cells.Item(42, 4).Value = 44529
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = 100114007
$ws.Cells.Item(42, 7).Value = "Jengibre"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 790
$ws.Cells.Item(42, 11).Value = 16000
$ws.Cells.Item(42, 12).Value = 18000
$ws.Cells.Item(42, 13).Value = 16987
$ws.Cells.Item(42, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(42, 15).Value = "Perú"
$ws.Cells.Item(42, 16).Value = 1307
$ws.Cells.Item(42, 17).Value = 13
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# New row 43: Jengibre, Segunda, fecha 2021-11-29 (serial 44529)
$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(43, 3).Value = "Metropolitana"
$ws.Cells.Item(43, 4).Value = 44529
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = 100114007
$ws.Cells.Item(43, 7).Value = "Jengibre"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Segunda"
$ws.Cells.Item(43, 10).Value = 430
$ws.Cells.Item(43, 11).Value = 13000
$ws.Cells.Item(43, 12).Value = 14000
$ws.Cells.Item(43, 13).Value = 13500
$ws.Cells.Item(43, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(43, 15).Value = "Perú"
$ws.Cells.Item(43, 16).Value = 1038
$ws.Cells.Item(43, 17).Value = 13
$ws.Cells.Item(43, 18).Value = "Hortaliza"
